# Fixed issue #13 Permitir que en los ficheros de metadatos dos columnas
# se puedan relacionar para crear SKOS jerarquicos.
#
# A new row is inserted right under the header row. It holds the
# "machine name" (slug) of each column header: lowercase, accents
# stripped, and any run of non [a-z0-9] characters collapsed to a
# single dash (leading/trailing dashes trimmed). This lets two columns
# be related to each other (e.g. a code column and its matching name
# column) when building the SKOS hierarchy from these metadata sheets.

function Slugify([string]$text) {
    $result = $text.ToLower()

    $accentMap = @{
        "á" = "a"; "à" = "a"; "ä" = "a"; "â" = "a"; "ã" = "a"; "å" = "a";
        "é" = "e"; "è" = "e"; "ë" = "e"; "ê" = "e";
        "í" = "i"; "ì" = "i"; "ï" = "i"; "î" = "i";
        "ó" = "o"; "ò" = "o"; "ö" = "o"; "ô" = "o"; "õ" = "o";
        "ú" = "u"; "ù" = "u"; "ü" = "u"; "û" = "u";
        "ñ" = "n"; "ç" = "c"
    }
    foreach ($accented in $accentMap.Keys) {
        $result = $result.Replace($accented, $accentMap[$accented])
    }

    $result = [System.Text.RegularExpressions.Regex]::Replace($result, "[^a-z0-9]+", "-")
    $result = $result.Trim("-")
    return $result
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedCols = $ws.UsedRange.Columns.Count

# Remember the header texts (row 1) before we shift rows around.
$headers = @{}
for ($col = 1; $col -le $usedCols; $col++) {
    $headers[$col] = $ws.Cells.Item(1, $col).Value2
}

# Insert a new row below the header row; this pushes the previous rows
# 2, 3, 4 (measure/dimension ids, "medida"/"dim" markers, xsd types /
# URI templates) down to rows 3, 4, 5 respectively.
$ws.Rows.Item(2).Insert()

# Fill the newly inserted row 2 with the slug of each column's header.
for ($col = 1; $col -le $usedCols; $col++) {
    $colLetter = [char](64 + $col)
    $slug = Slugify($headers[$col])
    $ws.Range("$colLetter" + "2").Value = $slug
}
